$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.233.29'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '1.862.16'
$ws.Range('E3').Value = '  -1.18%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7058'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.37'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07817'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3111'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.27'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08000'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.34%  '
$ws.Range('D12').Value = '1.864.93'
$ws.Range('E12').Value = '  -0.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '93.59'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.14%  '
$ws.Range('E14').Value = '  -1.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6951'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.354'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.18%  '
$ws.Range('D17').Value = '29.452.18'
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008269'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '251.87'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.42%  '
$ws.Range('D20').Value = '2.172.82'
$ws.Range('E20').Value = '  +1.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.550'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1555'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.988'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.47'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.47%  '
$ws.Range('E28').Value = '  +0.71%  '
$ws.Range('E29').Value = '  -0.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.268'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.266'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.211'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05269'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.888'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7446'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E36').Value = '  -2.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.707'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.21%  '
$ws.Range('E38').Value = '  -1.36%  '
$ws.Range('D39').Value = '1.249.92'
$ws.Range('E39').Value = '  -2.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.740'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.288'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9031'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '111.25'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.34%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '71.80'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.36%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '2.062.88'
$ws.Range('E46').Value = '  +2.07%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000128'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('E48').Value = '  -0.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.785'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.380'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.010'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.62%  '
